# -----------------------------------------------------------------------------
# Updates the cryptocurrency price list on Sheet1 with the latest scraped data.
# For each changed cell we write the literal display text exactly as produced by
# the scraper. Price values in column D are stored as *text* in the workbook
# (e.g. '41.816.60', '241.14', '0.622'); whenever the new text would otherwise be
# auto-recognised by Excel as a number (plain decimals like '241.14'), we prefix it
# with a leading apostrophe so Excel keeps storing/treating it as text, matching
# the original data's type.
# -----------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item('Sheet1')

# Row 2
$ws.Range('D2').Value = '41.816.60'
$ws.Range('E2').Value = '  -1.18%  '

# Row 3
$ws.Range('D3').Value = '2.219.97'
$ws.Range('E3').Value = '  -0.39%  '

# Row 4
$ws.Range('E4').Value = '  +0.03%  '

# Row 5
$ws.Range('D5').Value = '''241.14'
$ws.Range('E5').Value = '  -1.60%  '

# Row 6
$ws.Range('D6').Value = '''0.622'
$ws.Range('E6').Value = '  -1.15%  '

# Row 7
$ws.Range('D7').Value = '''72.73'
$ws.Range('E7').Value = '  -1.80%  '

# Row 8
$ws.Range('E8').Value = '  +0.15%  '

# Row 9
$ws.Range('D9').Value = '''0.589'
$ws.Range('E9').Value = '  -4.55%  '

# Row 10
$ws.Range('D10').Value = '''41.38'
$ws.Range('E10').Value = '  -3.56%  '

# Row 11
$ws.Range('D11').Value = '''0.0942'
$ws.Range('E11').Value = '  -2.28%  '

# Row 12
$ws.Range('E12').Value = '  -0.41%  '

# Row 13
$ws.Range('D13').Value = '''6.86'
$ws.Range('E13').Value = '  -4.02%  '

# Row 14
$ws.Range('D14').Value = '2.553.76'
$ws.Range('E14').Value = '  -0.26%  '

# Row 15
$ws.Range('D15').Value = '''14.15'
$ws.Range('E15').Value = '  -2.16%  '

# Row 16
$ws.Range('D16').Value = '''0.827'
$ws.Range('E16').Value = '  -2.84%  '

# Row 17
$ws.Range('D17').Value = '2.212.56'
$ws.Range('E17').Value = '  -0.89%  '

# Row 18
$ws.Range('D18').Value = '41.707.12'
$ws.Range('E18').Value = '  -1.09%  '

# Row 19
$ws.Range('D19').Value = '''0.0000104'
$ws.Range('E19').Value = '  -5.45%  '

# Row 20
$ws.Range('D20').Value = '''6.13'
$ws.Range('E20').Value = '  -0.55%  '

# Row 21
$ws.Range('D21').Value = '''71.70'
$ws.Range('E21').Value = '  -0.64%  '

# Row 22
$ws.Range('D22').Value = '''10.76'
$ws.Range('E22').Value = '  +7.45%  '

# Row 23
$ws.Range('D23').Value = '''228.28'
$ws.Range('E23').Value = '  -1.32%  '

# Row 24
$ws.Range('D24').Value = '''2.02'
$ws.Range('E24').Value = '  -6.34%  '

# Row 25
$ws.Range('E25').Value = '  +0.05%  '

# Row 26
$ws.Range('D26').Value = '''11.28'
$ws.Range('E26').Value = '  -4.69%  '

# Row 27
$ws.Range('E27').Value = '  -0.40%  '

# Row 28
$ws.Range('D28').Value = '''2.25'
$ws.Range('E28').Value = '  -2.12%  '

# Row 29
$ws.Range('D29').Value = '''2.21'
$ws.Range('E29').Value = '  -0.88%  '

# Row 30
$ws.Range('D30').Value = '''166.78'
$ws.Range('E30').Value = '  -0.27%  '

# Row 31
$ws.Range('D31').Value = '''20.36'
$ws.Range('E31').Value = '  -3.20%  '

# Row 32
$ws.Range('D32').Value = '''0.0789'
$ws.Range('E32').Value = '  -2.11%  '

# Row 33
$ws.Range('D33').Value = '''5.43'
$ws.Range('E33').Value = '  -5.72%  '

# Row 34
$ws.Range('D34').Value = '''30.50'
$ws.Range('E34').Value = '  +3.22%  '

# Row 35
$ws.Range('D35').Value = '''0.123'
$ws.Range('E35').Value = '  -1.68%  '

# Row 36
$ws.Range('E36').Value = '  -8.15%  '

# Row 37
$ws.Range('D37').Value = '''4.24'
$ws.Range('E37').Value = '  -4.27%  '

# Row 38
$ws.Range('D38').Value = '''0.0301'
$ws.Range('E38').Value = '  -1.88%  '

# Row 39
$ws.Range('D39').Value = '''12.95'
$ws.Range('E39').Value = '  -1.60%  '

# Row 40
$ws.Range('D40').Value = '''2.10'
$ws.Range('E40').Value = '  -2.94%  '

# Row 41
$ws.Range('B41').Value = 'THORChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D41').Value = '''5.60'
$ws.Range('E41').Value = '  -0.57%  '

# Row 42
$ws.Range('B42').Value = 'MultiversX'
$ws.Range('C42').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D42').Value = '''63.70'
$ws.Range('E42').Value = '  +0.91%  '

# Row 43
$ws.Range('D43').Value = '''0.195'
$ws.Range('E43').Value = '  -3.19%  '

# Row 44
$ws.Range('D44').Value = '''8.61'
$ws.Range('E44').Value = '  -2.54%  '

# Row 45
$ws.Range('D45').Value = '''101.61'
$ws.Range('E45').Value = '  -3.32%  '

# Row 46
$ws.Range('D46').Value = '''0.0991'
$ws.Range('E46').Value = '  -2.68%  '

# Row 47
$ws.Range('E47').Value = '  -1.06%  '

# Row 48
$ws.Range('D48').Value = '''1.16'
$ws.Range('E48').Value = '  -1.45%  '

# Row 49
$ws.Range('D49').Value = '''2.31'
$ws.Range('E49').Value = '  -3.35%  '

# Row 50
$ws.Range('E50').Value = '  -1.58%  '

# Row 51
$ws.Range('D51').Value = '2.428.54'
